$d = $word.ActiveDocument

# --- Change 1: remove the "in progress, not final" paragraph (bold+italic,
# FirstParagraph style) that immediately follows the "F2025" date paragraph
# and precedes the "overview" bookmark / Overview heading. ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "in progress, not final`r") {
        $p.Range.Delete()
        break
    }
}

# --- Change 2: prepend a bold "[in progress, not finalized]" note (followed
# by a plain space) in front of "Use the data at the" in the New York /
# election-data bullet (3rd item of the numId=1006 list). ---
$rng = $d.Content
$rng.Find.Execute("Use the data at the", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$insPos = $rng.Start
$noteText = "[in progress, not finalized]"

$ins = $d.Range($insPos, $insPos)
$ins.InsertBefore($noteText + " ")

$boldRng = $d.Range($insPos, $insPos + $noteText.Length)
$boldRng.Font.Bold = $true
